$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Target cluster is ECs (old rows 8, 5, 2),
# deleting bottom-to-top so row indices of earlier rows stay valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Refresh the remaining rows (now rows 2-7) with the recomputed TPM-based values.
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1510096666666667
$ws.Range("H2").Value = 0.453029
$ws.Range("I2").Value = 0.01105950042918124
$ws.Range("J2").Value = 0.01105950042918124
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.016376
$ws.Range("N2").Value = 0.049128
$ws.Range("O2").Value = 0.4917717717717718
$ws.Range("P2").Value = 0.4917717717717718
$ws.Range("Q2").Value = 0.002472934301333333
$ws.Range("R2").Value = 0.022256408712
$ws.Range("S2").Value = 0.005438750120969129
$ws.Range("T2").Value = 0.005438750120969129

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1510096666666667
$ws.Range("H3").Value = 0.453029
$ws.Range("I3").Value = 0.01105950042918124
$ws.Range("J3").Value = 0.01105950042918124
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.016924
$ws.Range("N3").Value = 0.050772
$ws.Range("O3").Value = 0.5082282282282282
$ws.Range("P3").Value = 0.5082282282282282
$ws.Range("Q3").Value = 0.002555687598666667
$ws.Range("R3").Value = 0.023001188388
$ws.Range("S3").Value = 0.005620750308212112
$ws.Range("T3").Value = 0.005620750308212112

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.23495333333333
$ws.Range("H4").Value = 30.70486
$ws.Range("I4").Value = 0.7495776481151314
$ws.Range("J4").Value = 0.7495776481151314
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016376
$ws.Range("N4").Value = 0.049128
$ws.Range("O4").Value = 0.4917717717717718
$ws.Range("P4").Value = 0.4917717717717718
$ws.Range("Q4").Value = 0.1676075957866666
$ws.Range("R4").Value = 1.50846836208
$ws.Range("S4").Value = 0.3686211280940959
$ws.Range("T4").Value = 0.3686211280940959

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.23495333333333
$ws.Range("H5").Value = 30.70486
$ws.Range("I5").Value = 0.7495776481151314
$ws.Range("J5").Value = 0.7495776481151314
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016924
$ws.Range("N5").Value = 0.050772
$ws.Range("O5").Value = 0.5082282282282282
$ws.Range("P5").Value = 0.5082282282282282
$ws.Range("Q5").Value = 0.1732163502133333
$ws.Range("R5").Value = 1.55894715192
$ws.Range("S5").Value = 0.3809565200210355
$ws.Range("T5").Value = 0.3809565200210355

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.268330666666666
$ws.Range("H6").Value = 9.804991999999999
$ws.Range("I6").Value = 0.2393628514556874
$ws.Range("J6").Value = 0.2393628514556874
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.016376
$ws.Range("N6").Value = 0.049128
$ws.Range("O6").Value = 0.4917717717717718
$ws.Range("P6").Value = 0.4917717717717718
$ws.Range("Q6").Value = 0.05352218299733331
$ws.Range("R6").Value = 0.4816996469759999
$ws.Range("S6").Value = 0.1177118935567068
$ws.Range("T6").Value = 0.1177118935567068

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.268330666666666
$ws.Range("H7").Value = 9.804991999999999
$ws.Range("I7").Value = 0.2393628514556874
$ws.Range("J7").Value = 0.2393628514556874
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016924
$ws.Range("N7").Value = 0.050772
$ws.Range("O7").Value = 0.5082282282282282
$ws.Range("P7").Value = 0.5082282282282282
$ws.Range("Q7").Value = 0.05531322820266665
$ws.Range("R7").Value = 0.4978190538239999
$ws.Range("S7").Value = 0.1216509578989806
$ws.Range("T7").Value = 0.1216509578989806

